$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.966.21'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '2.588.89'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.64%  '
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").Value = '2.588.74'
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '3.062.30'
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("E16").Value = '  -4.90%  '
$ws.Range("D17").Value = '66.974.88'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = '2.598.30'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '362.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("E21").Value = '  -4.61%  '
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '571.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("D31").Value = '0.0₃0972'
$ws.Range("E31").Value = '  -6.06%  '
$ws.Range("E32").Value = '  -4.66%  '
$ws.Range("E33").Value = '  -4.51%  '
$ws.Range("E34").Value = '  -4.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -6.42%  '
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").Value = '  -1.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("E44").Value = '  -4.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '150.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("D47").Value = '0.0₆0281'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0776'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E50").Value = '  -2.87%  '
$ws.Range("E51").Value = '  +0.79%  '
